# Avant-Projet / "Stalker In The Middle" - orthographe / mise en forme
#
# This script merges several two-run paragraphs (which only existed as two
# runs because of a leading-word / rest-of-sentence split) back into a
# single run, and fixes the spelling of "Etude" -> "Étude" on slide 9,
# matching the commit "evaluation bareme, Modification orthographe".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 16 - "Relations parties prenantes"
# "Envoi " + "ordre du jour, prise de notes, rédaction compte-rendu"
# -> "Envoi ordre du jour, prise de notes, rédaction compte-rendu"
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$tr16 = $s16.Shapes.Item(2).TextFrame.TextRange
$para = $tr16.Paragraphs(8)
$para.Runs(1).Text = "Envoi ordre du jour, prise de notes, rédaction compte-rendu"
$para.Runs(2).Text = ""

# ---------------------------------------------------------------------
# Slide 8 - "Fonctionnalités" (ARP spoofing)
# "Scanner " + "l'intégralité du réseau" -> "Scanner l'intégralité du réseau"
# "Modification " + "table ARP" -> "Modification table ARP"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange

$para = $tr8.Paragraphs(2)
$para.Runs(1).Text = "Scanner l’intégralité du réseau"
$para.Runs(2).Text = ""

$para = $tr8.Paragraphs(6)
$para.Runs(1).Text = "Modification table ARP"
$para.Runs(2).Text = ""

# ---------------------------------------------------------------------
# Slide 9 - "Fonctionnalités" (packet interception)
# "Modification et retransmission des " + "paquets"
# "Création " + "interface graphique"
# "Librairie " + "Curses"
# "Obtenir " + "mots de passe transitant en clair"
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange

$para = $tr9.Paragraphs(3)
$para.Runs(1).Text = "Modification et retransmission des paquets"
$para.Runs(2).Text = ""

$para = $tr9.Paragraphs(4)
$para.Runs(1).Text = "Création interface graphique"
$para.Runs(2).Text = ""

$para = $tr9.Paragraphs(5)
$para.Runs(1).Text = "Librairie Curses"
$para.Runs(2).Text = ""

$para = $tr9.Paragraphs(6)
$para.Runs(1).Text = "Obtenir mots de passe transitant en clair"
$para.Runs(2).Text = ""

# Spelling fix: "Etude des solutions pour intégrer les cookies sur Firefox"
# -> "Étude des solutions pour intégrer les cookies sur Firefox"
$para = $tr9.Paragraphs(9)
$para.Characters(1, 1).Text = "É"
